$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Drop the "episode 6..10" columns (H:L) entirely -- this both removes the
# cells and shrinks the sheet's used range/dimension down to column G.
$ws.Range("H1:L2").EntireColumn.Delete()

# Row 1: header label changes from "Episode" to "Values"; C1:G1 (1..5) stay.
$ws.Range("B1").Value = "Values"

# Row 2 used to hold the episode-1 results; it now becomes the
# "Final Value" metric row.
$ws.Range("B2").Value = "Final Value"
$ws.Range("C2").Value = 2775960.836417448
$ws.Range("D2").Value = 2160602.906874976
$ws.Range("E2").Value = 2347527.502617296
$ws.Range("F2").Value = 5485563.563985286
$ws.Range("G2").Value = 2914767.304096289

# Carry A2's formatting (bold/bordered/centered) down onto the two new
# index cells A3 and A4.
$ws.Range("A2").Copy()
$ws.Range("A3:A4").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Row 3: new "Annualized Return" metric row.
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "Annualized Return"
$ws.Range("C3").Value = 0.4010332452872412
$ws.Range("D3").Value = 0.2897389638376329
$ws.Range("E3").Value = 0.325572636022857
$ws.Range("F3").Value = 0.754472108596578
$ws.Range("G3").Value = 0.4237939922529597

# Row 4: new "Sharpe Ratio" metric row.
$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "Sharpe Ratio"
$ws.Range("C4").Value = 0.8370045488494169
$ws.Range("D4").Value = 0.5616484844575728
$ws.Range("E4").Value = 0.8522895704067394
$ws.Range("F4").Value = 1.632536598814653
$ws.Range("G4").Value = 1.043831200602436
